# Add Class: Floating Bus
#
# On the "Device" sheet, insert a new row above the existing "Notes:" row
# (row 16) to introduce a new bus-type class entry, "100, Floating Bus",
# styled to match the adjoining header-style cells (bold, same as A16).
#
# Also restores the various sheet selections that were left behind by the
# editing session on the other sheets (PowerFlow, NetworkLine) before the
# final selection on the Device sheet.

$wb = $excel.ActiveWorkbook

# --- PowerFlow sheet: leftover selection from browsing ---
$wsPowerFlow = $wb.Worksheets.Item("PowerFlow")
$wsPowerFlow.Range("F13").Select()

# --- NetworkLine sheet: leftover selection from browsing ---
$wsNetworkLine = $wb.Worksheets.Item("NetworkLine")
$wsNetworkLine.Range("E19").Select()

# --- Device sheet: the actual content edit ---
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Activate()

# Insert a new blank row at row 16, pushing "Notes:" (and everything below)
# down by one row.
$wsDevice.Rows.Item(16).Insert()

# Match the style of the neighboring label cell (A16, bold "header" style)
# and add the new class description text in column B.
$wsDevice.Range("A16").Copy($wsDevice.Range("B16"))
$wsDevice.Range("B16").Value = "100, Floating Bus"

# Leave the selection where the author ended up after the edit.
$wsDevice.Range("C25").Select()
